$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 39
$ws.Range("F2").Value = 21
$ws.Range("H2").Value = 23

$ws.Range("E3").Value = 35
$ws.Range("F3").Value = 20
$ws.Range("H3").Value = 23

$ws.Range("E6").Value = 65
$ws.Range("F6").Value = 32
$ws.Range("H6").Value = 39

$ws.Range("E8").Value = 52
$ws.Range("F8").Value = 21
$ws.Range("H8").Value = 27

$ws.Range("E9").Value = 19

$ws.Range("E11").Value = 24

$ws.Range("E12").Value = 43

$ws.Range("E14").Value = 45

$ws.Range("E15").Value = 126
$ws.Range("F15").Value = 62
$ws.Range("H15").Value = 73

$ws.Range("E16").Value = 351
$ws.Range("F16").Value = 110
$ws.Range("H16").Value = 198

$ws.Range("E17").Value = 39

$ws.Range("E18").Value = 108
$ws.Range("F18").Value = 36
$ws.Range("H18").Value = 59
